$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = -10.984
$ws.Range("C32").Value = -13.584
$ws.Range("C36").Value = -12.732
$ws.Range("C38").Value = -12.703
$ws.Range("C46").Value = -14.464
$ws.Range("C54").Value = -12.705
$ws.Range("C55").Value = -13.752
$ws.Range("C67").Value = -11.595
$ws.Range("C69").Value = -11.038
$ws.Range("C72").Value = -11.555
$ws.Range("C91").Value = -10.98
$ws.Range("C99").Value = -12.715
